$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently ends at row 21 (id 110020). We're appending 9 more
# user_detail records (ids 110021-110029). Duplicate the formatting of the
# last existing data row down into the new rows first (this is what gives
# the "email" column its left/no-fill style and the "is_active" boolean
# column its style, matching how the existing rows are formatted), then
# fill in the actual values.
$ws.Rows("21:21").Copy()
$ws.Rows("22:30").Insert(-4121)  # xlShiftDown, carries the copied formatting

$ids     = @(110021,110022,110023,110024,110025,110026,110027,110028,110029)
$uins    = @(7316931025,9137847236,8428758532,9804209494,7105248214,9316557128,8103486949,9601932866,9317596765)
# Note: the existing "name" column values in this workbook join first/last
# name with a U+00A0 (non-breaking space), not a plain space (consistent
# mock-data-generator artifact across all pre-existing rows) - match it.
$nbsp = [char]0x00A0
$names   = @("Magdalena${nbsp}Weber","Adrienne${nbsp}Hoffman","Adrienne${nbsp}Mcgee","Amare${nbsp}Coleman","Dawson${nbsp}Ibarra","Elvis${nbsp}Mcmillan","Steve${nbsp}George","Colton${nbsp}Elliott","Carolyn${nbsp}Rodriguez")
$emails  = @("magdalena.weber@xyz.com","adrienne.hoffman@xyz.com","adrienne.mcgee@xyz.com","amare.coleman@xyz.com","dawson.ibarra@xyz.com","elvis.mcmillan@xyz.com","steve.george@xyz.com","colton.elliott@xyz.com","carolyn.rodriguez@xyz.com")
$mobiles = @(932122450,848488000,894773246,956554588,765455583,884282274,971073663,809908673,818876429)

# Write column-by-column (not row-by-row) so that new shared-string entries
# land in the same order as the source workbook: all 9 names first, then
# all 9 emails.
for ($i = 0; $i -lt 9; $i++) {
    $ws.Range("A" + (22 + $i)).Value = $ids[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $ws.Range("B" + (22 + $i)).Value = $uins[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $ws.Range("C" + (22 + $i)).Value = $names[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $ws.Range("D" + (22 + $i)).Value = $emails[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $ws.Range("E" + (22 + $i)).Value = $mobiles[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $row = 22 + $i
    $ws.Range("F" + $row).Value = "ACT"
    $ws.Range("G" + $row).Value = "eng"
    $ws.Range("H" + $row).Value = "PWD"
    $ws.Range("I" + $row).Value = $true
    $ws.Range("J" + $row).Value = "superadmin"
    $ws.Range("K" + $row).Value = "now()"
}

# Match the author's final selection/scroll state.
$null = $ws.Range("A22:K30").Select()
